# Apply updated "Return_with_prediction" (G), "return_pct_change" (H),
# and "mean_return_pct_change" (I, row 2 only) values as produced by the
# new predicted prices / compared returns recalculation.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("G2").Value = 0.0383176576800623
$ws.Range("H2").Value = -20.5532567916236
$ws.Range("I2").Value = -52.57941970496287

$ws.Range("G3").Value = 0.04350316171216467
$ws.Range("H3").Value = 13.41765895898991

$ws.Range("G4").Value = -0.4791798042860266
$ws.Range("H4").Value = -5.833410596165616

$ws.Range("G5").Value = -0.4683252099896513
$ws.Range("H5").Value = 2.194862346172618

$ws.Range("G6").Value = 0.2487269103236668
$ws.Range("H6").Value = 6.461136346844026

$ws.Range("G7").Value = 0.242571720986055
$ws.Range("H7").Value = 9.972638828398464

$ws.Range("G8").Value = 0.1651096967438662
$ws.Range("H8").Value = -1.017191811472053

$ws.Range("G9").Value = 0.1679767773427478
$ws.Range("H9").Value = -2.345515506324745

$ws.Range("G10").Value = -0.01354062927564988
$ws.Range("H10").Value = -185.6985155209044

$ws.Range("G11").Value = -0.02635858095305886
$ws.Range("H11").Value = -79.85703968024008

$ws.Range("G12").Value = 0.1433146937789669
$ws.Range("H12").Value = 4.821759146787624

$ws.Range("G13").Value = 0.1435317737545024
$ws.Range("H13").Value = 15.15930530826543

$ws.Range("G14").Value = 0.255361814136216
$ws.Range("H14").Value = 3.245957106971121

$ws.Range("G15").Value = 0.2527791181297078
$ws.Range("H15").Value = 0.04041055342443028

$ws.Range("G16").Value = 0.1390085473030947
$ws.Range("H16").Value = -9.422553163762664

$ws.Range("G17").Value = 0.1434482756203169
$ws.Range("H17").Value = -5.010702805655138

$ws.Range("G18").Value = -0.01788526529322874
$ws.Range("H18").Value = -9.223619450182753

$ws.Range("G19").Value = -0.01076107448395934
$ws.Range("H19").Value = -1177.940345931511

$ws.Range("G20").Value = 0.1378603111761096
$ws.Range("H20").Value = -0.5782844001543764

$ws.Range("G21").Value = 0.1389899734116273
$ws.Range("H21").Value = -2.868385920880806

$ws.Range("G22").Value = 0.1763681344037855
$ws.Range("H22").Value = -5.289114936521217

$ws.Range("G23").Value = 0.1748823191194182
$ws.Range("H23").Value = -2.551082271397773

$ws.Range("G24").Value = -0.108462293233846
$ws.Range("H24").Value = -14.89982881706721

$ws.Range("G25").Value = -0.09145317539788259
$ws.Range("H25").Value = 8.187336312119392

$ws.Range("G26").Value = 0.2397016833863921
$ws.Range("H26").Value = 4.154524366191885

$ws.Range("G27").Value = 0.2511625436831053
$ws.Range("H27").Value = 7.993346750013923

$ws.Range("G28").Value = 0.04662143922384045
$ws.Range("H28").Value = -20.71623732721388

$ws.Range("G29").Value = 0.06771342341845395
$ws.Range("H29").Value = -4.067602832061471
